$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are purely numeric-looking text (e.g. "1.00", "0.0900")
# must be forced to Text format first so Excel does not coerce them to numbers
# and strip meaningful trailing/leading zeros or decimal formatting.
$textCells = @("D5", "D7", "D8", "D9", "D11", "D14", "D16", "D18", "D19", "D22", "D23", "D25", "D26", "D40", "D41", "D42", "D45", "D49", "D50")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = '29.891.77'
$ws.Range("E2").Value = '  -0.02%  '

# Row 3
$ws.Range("D3").Value = '1.635.55'
$ws.Range("E3").Value = '  +0.91%  '

# Row 4
$ws.Range("E4").Value = '  +0.78%  '

# Row 5
$ws.Range("D5").Value = '215.32'
$ws.Range("E5").Value = '  +0.90%  '

# Row 6
$ws.Range("E6").Value = '  +0.15%  '

# Row 7
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.83%  '

# Row 8
$ws.Range("D8").Value = '28.81'
$ws.Range("E8").Value = '  -1.20%  '

# Row 9
$ws.Range("D9").Value = '0.261'
$ws.Range("E9").Value = '  +0.55%  '

# Row 10
$ws.Range("E10").Value = '  +0.43%  '

# Row 11
$ws.Range("D11").Value = '0.0900'
$ws.Range("E11").Value = '  -1.21%  '

# Row 12
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.739.91'
$ws.Range("E12").Value = '  +7.04%  '

# Row 13
$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").Value = '1.871.35'
$ws.Range("E13").Value = '  +1.05%  '

# Row 14
$ws.Range("D14").Value = '0.589'
$ws.Range("E14").Value = '  +3.90%  '

# Row 15
$ws.Range("E15").Value = '  +6.87%  '

# Row 16
$ws.Range("D16").Value = '3.86'
$ws.Range("E16").Value = '  -1.28%  '

# Row 17
$ws.Range("D17").Value = '29.910.52'
$ws.Range("E17").Value = '  +0.13%  '

# Row 18
$ws.Range("D18").Value = '64.61'
$ws.Range("E18").Value = '  +0.34%  '

# Row 19
$ws.Range("D19").Value = '240.62'
$ws.Range("E19").Value = '  -0.64%  '

# Row 20
$ws.Range("E20").Value = '  -0.70%  '

# Row 21
$ws.Range("E21").Value = '  +0.61%  '

# Row 22
$ws.Range("D22").Value = '9.90'
$ws.Range("E22").Value = '  +3.21%  '

# Row 23
$ws.Range("D23").Value = '4.14'
$ws.Range("E23").Value = '  +0.95%  '

# Row 24
$ws.Range("E24").Value = '  +2.56%  '

# Row 25
$ws.Range("D25").Value = '157.31'
$ws.Range("E25").Value = '  +1.33%  '

# Row 26
$ws.Range("D26").Value = '15.52'
$ws.Range("E26").Value = '  -0.68%  '

# Row 27
$ws.Range("E27").Value = '  -0.58%  '

# Row 28
$ws.Range("E28").Value = '  +0.82%  '

# Row 29
$ws.Range("E29").Value = '  +0.66%  '

# Row 30
$ws.Range("E30").Value = '  +0.65%  '

# Row 31
$ws.Range("E31").Value = '  -0.55%  '

# Row 32
$ws.Range("E32").Value = '  +1.10%  '

# Row 33
$ws.Range("E33").Value = '  -0.63%  '

# Row 34
$ws.Range("D34").Value = '1.425.51'
$ws.Range("E34").Value = '  +0.69%  '

# Row 35
$ws.Range("E35").Value = '  +2.97%  '

# Row 36
$ws.Range("E36").Value = '  -0.78%  '

# Row 37
$ws.Range("E37").Value = '  -3.07%  '

# Row 38
$ws.Range("E38").Value = '  +1.52%  '

# Row 39
$ws.Range("E39").Value = '  +0.17%  '

# Row 40
$ws.Range("D40").Value = '76.33'
$ws.Range("E40").Value = '  +10.28%  '

# Row 41
$ws.Range("D41").Value = '0.560'
$ws.Range("E41").Value = '  +0.83%  '

# Row 42
$ws.Range("D42").Value = '0.833'
$ws.Range("E42").Value = '  +0.68%  '

# Row 43
$ws.Range("E43").Value = '  -0.18%  '

# Row 44
$ws.Range("E44").Value = '  -0.18%  '

# Row 45
$ws.Range("D45").Value = '1.00'
$ws.Range("E45").Value = '  +0.73%  '

# Row 46
$ws.Range("E46").Value = '  -1.72%  '

# Row 47
$ws.Range("E47").Value = '  -0.73%  '

# Row 48
$ws.Range("D48").Value = '1.779.62'
$ws.Range("E48").Value = '  +1.06%  '

# Row 49
$ws.Range("D49").Value = '48.83'
$ws.Range("E49").Value = '  -8.88%  '

# Row 50
$ws.Range("D50").Value = '92.81'
$ws.Range("E50").Value = '  +4.96%  '

# Row 51
$ws.Range("E51").Value = '  +0.78%  '
